# Update symbol list (crypto prices/volumes) as of the GitHub Actions run
# on Mon Dec 12 15:32:40 UTC 2022.
#
# The "Price" column (D) stores numeric-looking values as TEXT (the sheet
# is generated by a scraper and keeps exact decimal formatting, e.g.
# trailing zeros like "6.560"). A plain `.Value = "6.560"` assignment
# would let Excel auto-convert it to the number 6.56, losing the trailing
# zero, so each touched price cell is briefly formatted as Text ("@")
# while its value is written (so Excel stores it verbatim as a string)
# and then restored to the workbook's default "Normal" style so no
# formatting residue is left behind - only the cell content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2"  "276.51"
Set-TextValue "D3"  "20.91"
Set-TextValue "D4"  "6.218"
Set-TextValue "D6"  "3.578"
Set-TextValue "D7"  "6.560"
Set-TextValue "D8"  "1.483"
Set-TextValue "D9"  "0.8216"
Set-TextValue "D10" "0.01384"
Set-TextValue "D11" "0.1643"
Set-TextValue "D12" "0.08277"
Set-TextValue "D13" "0.03506"

Set-TextValue "D16" "3.771"
$ws.Range("E16").Value = "15MCDexMCBBestin24h"

Set-TextValue "D17" "0.001611"
Set-TextValue "D18" "0.04702"
Set-TextValue "D19" "0.006422"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.807"
Set-TextValue "D24" "2.283"
Set-TextValue "D25" "0.3385"
Set-TextValue "D26" "0.1201"
Set-TextValue "D40" "0.04667"
Set-TextValue "D41" "0.007016"

# Rows 42 and 43 swapped identity (BKEXToken/CEJI reordered) with new data
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1105"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003522"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D45" "0.00006220"
Set-TextValue "D48" "0.001969"
